# Update the "Phieu van dap" submission sheet:
# - Row 45 (3.3 Phat tan tin rao vat): completion ratio 1 -> 0.75
# - Row 50 (4.2 Quan ly chuyen muc chinh, chuyen muc con): completion ratio 0.5 -> 0.1
# Downstream formulas (G42, G46, G11, H11) recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PhieuVanDap")

$ws.Range("D45").Value = 0.75
$ws.Range("D50").Value = 0.1

# Reset the view: scroll back to the top-left and clear any stray selection
# that was left over from editing further down the sheet.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A1").Select()
